$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.256247043609619
$ws.Range("B1").Value = 2.251749515533447
$ws.Range("C1").Value = 4.539510250091553
$ws.Range("D1").Value = 2.918430805206299
$ws.Range("E1").Value = 1.370581984519958
